$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 639 (shifts existing rows 639-650 down to 643-654)
$ws.Range("A639:A642").EntireRow.Insert()

# Common/static values for these new rows (same as neighboring rows in this block)
$mercadoId = 9
$mercado = "Vega Central Mapocho de Santiago"
$region = "Metropolitana"
$fecha = 44448
$codreg = 13
$tipo = "Fruta"
$productoId = 100108
$producto = "Tropicales y subtropicales"
$categoriaId = 100108005
$categoria = "Piña"
$variedad = "Caramelo"
$origen = "Ecuador"

$rows = @(
    @{ Row=639; Calidad="Especial"; Volumen=30; PMin=22000; PMax=22000; PProm=22000; Unidad="`$/caja 10 unidades"; PKg=2200; KgUnidad=10 },
    @{ Row=640; Calidad="Primera";  Volumen=25; PMin=22000; PMax=22000; PProm=22000; Unidad="`$/caja 12 unidades"; PKg=1833; KgUnidad=12 },
    @{ Row=641; Calidad="Segunda";  Volumen=30; PMin=22000; PMax=22000; PProm=22000; Unidad="`$/caja 14 unidades"; PKg=1571; KgUnidad=14 },
    @{ Row=642; Calidad="Tercera";  Volumen=20; PMin=22000; PMax=22000; PProm=22000; Unidad="`$/caja 16 unidades"; PKg=1375; KgUnidad=16 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
